# Add Info to Philadelphia (9) -- column K (header "US9")
# Rows where column K goes from 0 -> 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$rows = @(12, 15, 16, 17, 18, 20, 21, 22, 26, 27, 38, 40, 41)
foreach ($r in $rows) {
    $ws.Range("K$r").Value = 1
}

# Match the final active selection recorded in the workbook (K38).
$ws.Range("K38").Select()
